$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC-2020")
$ws.Activate()

$ws.Range("A1").Value = "test"
